$d = $word.ActiveDocument

$pairs = @(
    @("22×43=", "59×17="),
    @("93×30=", "26×31="),
    @("20×23=", "80×89="),
    @("99×17=", "97×17="),
    @("52×23=", "52×68="),
    @("37×41=", "29×49="),
    @("40×52=", "57×30="),
    @("82×48=", "91×24="),
    @("92×15=", "57×33="),
    @("15×16=", "84×15="),
    @("49×82=", "86×48="),
    @("57×93=", "74×91="),
    @("97×89=", "73×28="),
    @("54×39=", "94×69="),
    @("36×20=", "83×44="),
    @("69×70=", "97×71="),
    @("63×17=", "57×41="),
    @("54×73=", "44×72="),
    @("35×76=", "84×69="),
    @("24×88=", "87×55="),
    @("64×88=", "93×80="),
    @("71×97=", "67×54="),
    @("54×18=", "44×37="),
    @("24×78=", "38×17="),
    @("84×81=", "26×19=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
